$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row 2 / row 1(H-J duplicate table): "# sources cited / 10" ->
#    "# sources cited" (values below are no longer divided by ten).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "# sources cited"
$ws.Range("H2").Value = "# sources cited"

# ---------------------------------------------------------------------------
# 2. Quality-measure values for the Wikipedia side (columns C & F) are no
#    longer divided by ten -- multiply the old numbers by ten.
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = 128
$ws.Range("F3").Value = 208

$ws.Range("C4").Value = 139
$ws.Range("F4").Value = 179

$ws.Range("C5").Value = 39
$ws.Range("F5").Value = 59

$ws.Range("C6").Value = 7
$ws.Range("F6").Value = 77

$ws.Range("C7").Value = 35
$ws.Range("F7").Value = 45

$ws.Range("C8").Value = 14
$ws.Range("F8").Value = 14

# ---------------------------------------------------------------------------
# 3. Replace the comparison encyclopedic sources (column G) with the new
#    FactMonster based sources, and refresh their quality-measure values.
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = "Migration of Animals - FactoMonster"
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = -5
$ws.Range("J3").Value = 0

$ws.Range("G4").Value = "Tea - FactMonster"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = -4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = -2

$ws.Range("G5").Value = "Hibernation - FactMonster"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = -3
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -2

$ws.Range("G6").Value = "Electricity - FactMonster"
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = -3
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 2

$ws.Range("G7").Value = "n/a"
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0

$ws.Range("G8").Value = "Mummy - FactMonster"
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = -3
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3

# ---------------------------------------------------------------------------
# 4. Hyperlinks: rebuild the hyperlink collection. Clearing via the
#    Hyperlinks collection wipes every hyperlink on the sheet, so do that
#    once and then re-add every link we still want -- the wiki links in
#    column B are unchanged, the column G links now point at FactMonster
#    pages (G7 becomes "n/a" and keeps no hyperlink at all).
# ---------------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), "https://en.wikipedia.org/wiki/Bird_migration", ":~:text=Bird%20migration%20is%20the%20regular,Many%20species%20of%20bird%20migrate.&text=It%20occurs%20mainly%20in%20the,Sea%20or%20the%20Caribbean%20Sea.")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://en.wikipedia.org/wiki/Tea")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://en.wikipedia.org/wiki/Hibernation")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://en.wikipedia.org/wiki/Electricity")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://en.wikipedia.org/wiki/Water_cycle")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://en.wikipedia.org/wiki/Mummy")

$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.factmonster.com/dk/encyclopedia/animals/migration")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.factmonster.com/encyclopedia/plants/applied/food/tea")
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.factmonster.com/encyclopedia/science/biology/general/hibernation")
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.factmonster.com/math-science/physics/electricity")
$ws.Hyperlinks.Add($ws.Range("G8"), "https://www.factmonster.com/encyclopedia/people/modern/archaeology/mummy")

# ---------------------------------------------------------------------------
# 5. Methodology / Dictionary table (rows 14-18): update the wording of the
#    "# of sources cited" and "# of further reading sources" descriptions,
#    merge B:D on every row, center + wrap the merged cells, and resize the
#    rows.
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = "Examining the web page and counting the number of sources cited. Either in explicit citations or written as an entry in a bibliogrpahy."
$ws.Range("B17").Value = "Examining the web page and counting the number of further reading links cited, if present."

$ws.Range("B14:D14").Merge()
$ws.Range("B15:D15").Merge()
$ws.Range("B16:D16").Merge()
$ws.Range("B17:D17").Merge()
$ws.Range("B18:D18").Merge()

$rng1418 = $ws.Range("B14:D18")
$rng1418.HorizontalAlignment = -4108
$rng1418.WrapText = $true

$ws.Rows.Item(14).RowHeight = 38
$ws.Rows.Item(15).RowHeight = 69
$ws.Rows.Item(16).RowHeight = 50
$ws.Rows.Item(17).RowHeight = 51
$ws.Rows.Item(18).RowHeight = 34

# ---------------------------------------------------------------------------
# 6. Selection cosmetic change.
# ---------------------------------------------------------------------------
$ws.Range("B3").Select()

Write-Output "Workbook updated"
